$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.231077503630729
$ws.Range("C2").Value = 0.342661379653407
$ws.Range("E2").Value = 0.695969353980658
$ws.Range("F2").Value = 2.188743752544553
$ws.Range("G2").Value = 0.4075221873743686
$ws.Range("H2").Value = 0.5499249562740687
$ws.Range("I2").Value = 0.3491459723963111
$ws.Range("J2").Value = 0.03620267388696519

$ws.Range("B3").Value = 1.085047152610059
$ws.Range("C3").Value = 0.2989601208352042
$ws.Range("E3").Value = 0.6711336273832842
$ws.Range("F3").Value = 2.15598724742614
$ws.Range("G3").Value = 0.4096301497465902
$ws.Range("H3").Value = 0.5578365375708643
$ws.Range("I3").Value = 0.3600529139786293
$ws.Range("J3").Value = 0.03662848168936605

$ws.Range("B4").Value = 0.9952412421958456
$ws.Range("C4").Value = 0.2720566125698554
$ws.Range("E4").Value = 0.6561075010668844
$ws.Range("F4").Value = 2.137531834642672
$ws.Range("G4").Value = 0.4116184581827298
$ws.Range("H4").Value = 0.5632432898832178
$ws.Range("I4").Value = 0.367248341222929
$ws.Range("J4").Value = 0.03694118266878377

$ws.Range("B5").Value = 0.9586099889535262
$ws.Range("C5").Value = 0.2610755948907695
$ws.Range("E5").Value = 0.6500404276958989
$ws.Range("F5").Value = 2.130425970492453
$ws.Range("G5").Value = 0.4126017383420191
$ws.Range("H5").Value = 0.5655840729329498
$ws.Range("I5").Value = 0.370305215577087
$ws.Range("J5").Value = 0.03708140620129896

$ws.Range("B6").Value = 0.9525253446261672
$ws.Range("C6").Value = 0.2592511464545453
$ws.Range("E6").Value = 0.6490363917176438
$ws.Range("F6").Value = 2.129271050130995
$ws.Range("G6").Value = 0.4127754182392493
$ws.Range("H6").Value = 0.5659810455524905
$ws.Range("I6").Value = 0.370820317721634
$ws.Range("J6").Value = 0.03710546030970363

$ws.Range("B7").Value = 0.9947473587563991
$ws.Range("C7").Value = 0.2719085898222602
$ws.Range("E7").Value = 0.6560254507335372
$ws.Range("F7").Value = 2.137434325155581
$ws.Range("G7").Value = 0.4116310204304909
$ws.Range("H7").Value = 0.563274302575195
$ws.Range("I7").Value = 0.3672890633800794
$ws.Range("J7").Value = 0.03694302208700506

$ws.Range("B8").Value = 1.180756554701929
$ws.Range("C8").Value = 0.3276079591793462
$ws.Range("E8").Value = 0.6873597334094939
$ws.Range("F8").Value = 2.177104219702812
$ws.Range("G8").Value = 0.4081040753879392
$ws.Range("H8").Value = 0.5525385824293068
$ws.Range("I8").Value = 0.3528028112829986
$ws.Range("J8").Value = 0.0363387989458559

$ws.Range("B9").Value = 1.544352055034096
$ws.Range("C9").Value = 0.4362710445480502
$ws.Range("E9").Value = 0.7505764052701664
$ws.Range("F9").Value = 2.268137156983641
$ws.Range("G9").Value = 0.4067618475223469
$ws.Range("H9").Value = 0.535867525018574
$ws.Range("I9").Value = 0.3283808565932738
$ws.Range("J9").Value = 0.03556478952012654

$ws.Range("B10").Value = 1.81075456907098
$ws.Range("C10").Value = 0.5157687319832007
$ws.Range("E10").Value = 0.7981067797691281
$ws.Range("F10").Value = 2.343227853896252
$ws.Range("G10").Value = 0.4092677456574734
$ws.Range("H10").Value = 0.5263258236694526
$ws.Range("I10").Value = 0.3129084447426571
$ws.Range("J10").Value = 0.03525243979321147

$ws.Range("B11").Value = 1.931787634703028
$ws.Range("C11").Value = 0.5518631503290408
$ws.Range("E11").Value = 0.8199669523956743
$ws.Range("F11").Value = 2.379201470494053
$ws.Range("G11").Value = 0.4111862861479096
$ws.Range("H11").Value = 0.5225804057925387
$ws.Range("I11").Value = 0.3064148175467913
$ws.Range("J11").Value = 0.03516726407884008

$ws.Range("B12").Value = 1.977596762071698
$ws.Range("C12").Value = 0.5655211796780577
$ws.Range("E12").Value = 0.8282791423926739
$ws.Range("F12").Value = 2.393086949422468
$ws.Range("G12").Value = 0.4120263886364341
$ws.Range("H12").Value = 0.5212483271469921
$ws.Range("I12").Value = 0.3040349371639728
$ws.Range("J12").Value = 0.03514329762004209

$ws.Range("B13").Value = 1.967732009438862
$ws.Range("C13").Value = 0.5625801321243671
$ws.Range("E13").Value = 0.8264874426052273
$ws.Range("F13").Value = 2.390084729670605
$ws.Range("G13").Value = 0.4118403794272041
$ws.Range("H13").Value = 0.5215313686713614
$ws.Range("I13").Value = 0.3045439560746104
$ws.Range("J13").Value = 0.03514808897170241

$ws.Range("B14").Value = 1.935556856300423
$ws.Range("C14").Value = 0.5529870087474364
$ws.Range("E14").Value = 0.8206501158216923
$ws.Range("F14").Value = 2.380338552988974
$ws.Range("G14").Value = 0.4112531159519932
$ws.Range("H14").Value = 0.5224690829845855
$ws.Range("I14").Value = 0.3062174333665588
$ws.Range("J14").Value = 0.03516512565021657

$ws.Range("B15").Value = 1.915845554528346
$ws.Range("C15").Value = 0.5471096183711666
$ws.Range("E15").Value = 0.817079039053425
$ws.Range("F15").Value = 2.374403058960269
$ws.Range("G15").Value = 0.4109082417665206
$ws.Range("H15").Value = 0.5230547090258142
$ws.Range("I15").Value = 0.3072528124727754
$ws.Range("J15").Value = 0.03517664349925909

$ws.Range("B16").Value = 1.802841554608563
$ws.Range("C16").Value = 0.5134084654556546
$ws.Range("E16").Value = 0.7966829573283007
$ws.Range("F16").Value = 2.340913587137322
$ws.Range("G16").Value = 0.4091581712368679
$ws.Range("H16").Value = 0.5265826267995521
$ws.Range("I16").Value = 0.3133438477916535
$ws.Range("J16").Value = 0.03525916032144494

$ws.Range("B17").Value = 1.733476744073243
$ws.Range("C17").Value = 0.4927159922639817
$ws.Range("E17").Value = 0.7842316281434734
$ws.Range("F17").Value = 2.320835122806756
$ws.Range("G17").Value = 0.4082851361546886
$ws.Range("H17").Value = 0.5288997823838599
$ws.Range("I17").Value = 0.3172205944138611
$ws.Range("J17").Value = 0.03532443225448034

$ws.Range("B18").Value = 1.69356544132421
$ws.Range("C18").Value = 0.4808076977041651
$ws.Range("E18").Value = 0.7770923922070807
$ws.Range("F18").Value = 2.309457180213741
$ws.Range("G18").Value = 0.4078561767031346
$ws.Range("H18").Value = 0.5302885284232275
$ws.Range("I18").Value = 0.3195016344724309
$ws.Range("J18").Value = 0.03536732571173928

$ws.Range("B19").Value = 1.680049707492742
$ws.Range("C19").Value = 0.4767746368762573
$ws.Range("E19").Value = 0.7746790255134641
$ws.Range("F19").Value = 2.305634050599537
$ws.Range("G19").Value = 0.4077234587243623
$ws.Range("H19").Value = 0.5307683287870049
$ws.Range("I19").Value = 0.3202827322934034
$ws.Range("J19").Value = 0.03538276454956346

$ws.Range("B20").Value = 1.740862255989498
$ws.Range("C20").Value = 0.4949194184545718
$ws.Range("E20").Value = 0.7855547714673747
$ws.Range("F20").Value = 2.322954831718732
$ws.Range("G20").Value = 0.4083704858378212
$ws.Range("H20").Value = 0.5286473188915579
$ws.Range("I20").Value = 0.3168025996580397
$ws.Range("J20").Value = 0.03531692941257347

$ws.Range("B21").Value = 1.945008125069535
$ws.Range("C21").Value = 0.5558050180727037
$ws.Range("E21").Value = 0.8223637519919862
$ws.Range("F21").Value = 2.383194085115917
$ws.Range("G21").Value = 0.4114225135648581
$ws.Range("H21").Value = 0.5221913081284697
$ws.Range("I21").Value = 0.3057237392305829
$ws.Range("J21").Value = 0.03515989580933976

$ws.Range("B22").Value = 2.078291931119907
$ws.Range("C22").Value = 0.5955380905508036
$ws.Range("E22").Value = 0.8466199699068255
$ws.Range("F22").Value = 2.424098186202713
$ws.Range("G22").Value = 0.4140800567677303
$ws.Range("H22").Value = 0.5184748481114383
$ws.Range("I22").Value = 0.2989445658206549
$ws.Range("J22").Value = 0.03510562099303627

$ws.Range("B23").Value = 2.007168809591008
$ws.Range("C23").Value = 0.5743372648715876
$ws.Range("E23").Value = 0.8336557404553702
$ws.Range("F23").Value = 2.402125802672884
$ws.Range("G23").Value = 0.412600483222036
$ws.Range("H23").Value = 0.520412161033164
$ws.Range("I23").Value = 0.302520256765451
$ws.Range("J23").Value = 0.0351301295630968

$ws.Range("B24").Value = 1.737523367349752
$ws.Range("C24").Value = 0.4939232865831968
$ws.Range("E24").Value = 0.784956518566446
$ws.Range("F24").Value = 2.321995996171751
$ws.Range("G24").Value = 0.4083316721079484
$ws.Range("H24").Value = 0.5287612814349245
$ws.Range("I24").Value = 0.3169914124007036
$ws.Range("J24").Value = 0.03532030473605374

$ws.Range("B25").Value = 1.446117009894238
$ws.Range("C25").Value = 0.406934493105723
$ws.Range("E25").Value = 0.733284523438229
$ws.Range("F25").Value = 2.242078757402069
$ws.Range("G25").Value = 0.4065185160187639
$ws.Range("H25").Value = 0.5399047566750994
$ws.Range("I25").Value = 0.3345567646625582
$ws.Range("J25").Value = 0.03572959130970688
